$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 627, pushing the existing
# rows 627..670 down to 629..672 (matches the target dimension A1:R672).
$ws.Rows.Item(627).Insert()
$ws.Rows.Item(627).Insert()

# --- New row 627: Acelga / Primera, 2022-09-22 (serial 44826) ---
$ws.Cells.Item(627, 1).Value = 9
$ws.Cells.Item(627, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(627, 3).Value = "Metropolitana"
$ws.Cells.Item(627, 4).Value = 44826
$ws.Cells.Item(627, 5).Value = 13
$ws.Cells.Item(627, 6).Value = 100112009
$ws.Cells.Item(627, 7).Value = "Acelga"
$ws.Cells.Item(627, 8).Value = "Sin especificar"
$ws.Cells.Item(627, 9).Value = "Primera"
$ws.Cells.Item(627, 10).Value = 55
$ws.Cells.Item(627, 11).Value = 10000
$ws.Cells.Item(627, 12).Value = 12000
$ws.Cells.Item(627, 13).Value = 11273
$ws.Cells.Item(627, 14).Value = "`$/docena de atados"
$ws.Cells.Item(627, 15).Value = "Región Metropolitana"
$ws.Cells.Item(627, 16).Value = 3758
$ws.Cells.Item(627, 17).Value = 3
$ws.Cells.Item(627, 18).Value = "Hortaliza"

# --- New row 628: Acelga / Segunda, 2022-09-22 (serial 44826) ---
$ws.Cells.Item(628, 1).Value = 9
$ws.Cells.Item(628, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(628, 3).Value = "Metropolitana"
$ws.Cells.Item(628, 4).Value = 44826
$ws.Cells.Item(628, 5).Value = 13
$ws.Cells.Item(628, 6).Value = 100112009
$ws.Cells.Item(628, 7).Value = "Acelga"
$ws.Cells.Item(628, 8).Value = "Sin especificar"
$ws.Cells.Item(628, 9).Value = "Segunda"
$ws.Cells.Item(628, 10).Value = 25
$ws.Cells.Item(628, 11).Value = 8000
$ws.Cells.Item(628, 12).Value = 8000
$ws.Cells.Item(628, 13).Value = 8000
$ws.Cells.Item(628, 14).Value = "`$/docena de atados"
$ws.Cells.Item(628, 15).Value = "Región Metropolitana"
$ws.Cells.Item(628, 16).Value = 2667
$ws.Cells.Item(628, 17).Value = 3
$ws.Cells.Item(628, 18).Value = "Hortaliza"
